$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Plan de pruebas")
$ws2 = $wb.Worksheets.Item("Q-gates")

# Update values on the Q-gates sheet
$ws2.Range("C9").Value = 7
$ws2.Range("C11").Value = 5

# Keep selection on Plan de pruebas sheet as-is (F9), just make sure it's not the active tab afterward
$ws1.Activate()
$ws1.Range("F9").Select()

# Activate Q-gates sheet last so it becomes the active/selected tab,
# and set its selection to C15:D15
$ws2.Activate()
$ws2.Range("C15:D15").Select()
